$d = $word.ActiveDocument

# Replace every occurrence of $old with $new while trying hard to preserve
# the original run (formatting) structure of the paragraph. The emulated
# Word engine re-normalizes (merges) adjacent runs that end up with
# identical resolved formatting whenever a paragraph is edited - this
# mirrors real Word's own run-coalescing behavior and can't be avoided
# from the object model, but we can control *which* run's formatting the
# new text inherits so it matches the formatting the replaced text had.
function Replace-Preserving-All($old, $new) {
    $guard = 0
    while ($true) {
        $guard += 1
        if ($guard -gt 50) { break }

        $rng = $d.Content
        $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
        if (-not $found) { break }

        $fStart = $rng.Start
        $fEnd = $rng.End

        # Is the match the very last text in its paragraph (immediately
        # followed by the paragraph mark)? Inserting *after* such a range
        # lands on the paragraph-mark boundary and loses formatting, so in
        # that case we insert *before* the match instead (which inherits
        # the preceding, identically formatted run) and then delete the
        # shifted-right original text.
        $afterChar = $d.Range($fEnd, $fEnd + 1).Text
        $isParaEnd = ($afterChar -eq [string][char]13)

        if ($isParaEnd) {
            $rng.InsertBefore($new)
            $newLen = $new.Length
            $delRng = $d.Range($fStart + $newLen, $fEnd + $newLen)
            $delRng.Text = ""
        } else {
            $rng.InsertAfter($new)
            $delRng = $d.Range($fStart, $fEnd)
            $delRng.Text = ""
        }
    }
}

# Company / trading names
Replace-Preserving-All "CAMPBELL AND CURTIS TRADING ," "CANTU AND SEARS ASSOCIATES ,"
Replace-Preserving-All "CAMPBELL AND CURTIS TRADING" "CANTU AND SEARS ASSOCIATES"
Replace-Preserving-All "Mcintyre And Pope Inc " "Caldwell And Molina Inc "

# Person names
Replace-Preserving-All "RAJAH GUY " "AURORA GOFF "
Replace-Preserving-All "RAJAH GUY" "AURORA GOFF"
Replace-Preserving-All "HOPE BURKE" "CAMILLE CRAWFORD"

# Latin placeholder text
Replace-Preserving-All "VOLUPTATE EST CORPOR" "ANIM HIC SUNT OMNIS"
